$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value (values must stay as text, matching the
# original inline-string cell type, not be auto-converted to number/percentage)
$updates = @{
    'D2' = '304.47'
    'E2' = '5.92%'
    'E3' = '8.95%'
    'D4' = '5.318'
    'E4' = '4.25%'
    'E5' = '7.23%'
    'D6' = '7.812'
    'E6' = '5.39%'
    'D7' = '3.857'
    'E7' = '7.77%'
    'D8' = '1.466'
    'E8' = '6.11%'
    'D9' = '0.9235'
    'E9' = '2.32%'
    'D10' = '0.01752'
    'E10' = '2,602.71%'
    'D11' = '0.1699'
    'E11' = '6.50%'
    'D12' = '0.07693'
    'E12' = '8.17%'
    'D13' = '0.08068'
    'E13' = '5.12%'
    'D14' = '0.03065'
    'E14' = '5.27%'
    'D15' = '0.09887'
    'E15' = '10.00%'
    'D16' = '0.001497'
    'E16' = '-6.18%'
    'D17' = '0.04573'
    'E17' = '1.03%'
    'D18' = '0.006464'
    'E18' = '0.95%'
    'D19' = '3.479'
    'E19' = '0.26%'
    'E20' = '-0.06%'
    'D21' = '0.3336'
    'D22' = '0.1347'
    'D23' = '4.516'
    'E23' = '12.75%'
    'D24' = '0.1625'
    'E24' = '4.60%'
    'E25' = '1.25%'
    'D26' = '0.004425'
    'E26' = '0.94%'
    'D27' = '0.0001401'
    'E27' = '20.16%'
    'D28' = '0.0001742'
    'E28' = '8.11%'
    'D40' = '0.04537'
    'E40' = '5.87%'
    'D41' = '0.007203'
    'E41' = '5.71%'
    'E42' = '7.55%'
    'D43' = '0.002212'
    'E43' = '0.87%'
    'D44' = '0.01272'
    'E44' = '10.43%'
    'D45' = '0.00006168'
    'E45' = '7.38%'
    'D46' = '1.872'
    'E46' = '-2.95%'
    'D47' = '0.01301'
    'E47' = '-0.07%'
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    # Force text storage so strings like "304.47" or "5.92%" are not
    # reinterpreted by Excel as numbers/percentages.
    $range.NumberFormat = "@"
    $range.Value = $updates[$cell]
}
